$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 125, shifting existing rows 125:252 down to 126:253
$ws.Rows("125:125").Insert()

# Populate the newly inserted row 125 with the new data record
$ws.Range("A125").Value = 8
$ws.Range("B125").Value = "Terminal La Palmera de La Serena"
$ws.Range("C125").Value = "Coquimbo"
$ws.Range("D125").Value = 44874
$ws.Range("E125").Value = 4
$ws.Range("F125").Value = 100112037
$ws.Range("G125").Value = "Cebollín"
$ws.Range("H125").Value = "Sin especificar"
$ws.Range("I125").Value = "Primera"
$ws.Range("J125").Value = 1660
$ws.Range("K125").Value = 1200
$ws.Range("L125").Value = 1400
$ws.Range("M125").Value = 1300
$ws.Range("N125").Value = "$/paquete 6 unidades"
$ws.Range("O125").Value = "Provincia del Elquí"
$ws.Range("P125").Value = 217
$ws.Range("Q125").Value = 6
$ws.Range("R125").Value = "Hortaliza"
